$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 669.9231
$ws.Range("I32").Value = 733
$ws.Range("J32").Value = 651
$ws.Range("K32").Value = 733
$ws.Range("L32").Value = 651
$ws.Range("M32").Value = -407
$ws.Range("N32").Value = -1303

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3917.3157
$ws.Range("I100").Value = 2963.2222
$ws.Range("K100").Value = 2963.2222
$ws.Range("M100").Value = -2422.2222

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 10870523
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 10870523
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 32611569
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -32613785

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 14794
$ws.Range("I43").Value = 14499
$ws.Range("J43").Value = 14892.333
$ws.Range("K43").Value = 14499
$ws.Range("L43").Value = 14892.333
$ws.Range("N43").Value = -15518.333
$ws.Range("M43").Value = -14186

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1456.2333
$ws.Range("I45").Value = 1121
$ws.Range("J45").Value = 2378.125
$ws.Range("K45").Value = 1121
$ws.Range("L45").Value = 2378.125
$ws.Range("M45").Value = -744
$ws.Range("N45").Value = -3132.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 5966.6665
$ws.Range("I46").Value = 4000
$ws.Range("J46").Value = 6950
$ws.Range("K46").Value = 4000
$ws.Range("L46").Value = 6950
$ws.Range("N46").Value = -7588
$ws.Range("M46").Value = -3681

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 33333.332
$ws.Range("J123").Value = 33333.332
$ws.Range("L123").Value = 33333.332
$ws.Range("N123").Value = -43133.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2488.0417
$ws.Range("I134").Value = 1945.8422
$ws.Range("J134").Value = 4548.4
$ws.Range("K134").Value = 5837.5266
$ws.Range("L134").Value = 13645.2
$ws.Range("M134").Value = -3302.5266
$ws.Range("N134").Value = -18715.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1854212.5
$ws.Range("I31").Value = 2441176.2
$ws.Range("J31").Value = 3019.6924
$ws.Range("K31").Value = 2441176.2
$ws.Range("L31").Value = 3019.6924
$ws.Range("M31").Value = -2440881.2
$ws.Range("N31").Value = -3609.6924

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1854212.5
$ws.Range("I34").Value = 2441176.2
$ws.Range("J34").Value = 3019.6924
$ws.Range("K34").Value = 2441176.2
$ws.Range("L34").Value = 3019.6924
$ws.Range("M34").Value = -2440974.2
$ws.Range("N34").Value = -3423.6924

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 11366644
$ws.Range("I58").Value = 2471.8
$ws.Range("J58").Value = 26319502
$ws.Range("K58").Value = 2471.8
$ws.Range("L58").Value = 26319502
$ws.Range("M58").Value = -2268.8
$ws.Range("N58").Value = -26319908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 27282
$ws.Range("J88").Value = 27282
$ws.Range("L88").Value = 27282
$ws.Range("N88").Value = -28094

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 27282
$ws.Range("J91").Value = 27282
$ws.Range("L91").Value = 27282
$ws.Range("N91").Value = -30090

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 11366644
$ws.Range("I136").Value = 2471.8
$ws.Range("J136").Value = 26319502
$ws.Range("K136").Value = 7415.400000000001
$ws.Range("L136").Value = 78958506
$ws.Range("M136").Value = -4865.400000000001
$ws.Range("N136").Value = -78963606

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 9997.691999999999
$ws.Range("J64").Value = 15125
$ws.Range("L64").Value = 45375
$ws.Range("N64").Value = -45915

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 9997.691999999999
$ws.Range("J67").Value = 15125
$ws.Range("L67").Value = 45375
$ws.Range("N67").Value = -47247

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4367.1113
$ws.Range("I132").Value = 5401.1665
$ws.Range("J132").Value = 3539.8667
$ws.Range("K132").Value = 16203.4995
$ws.Range("L132").Value = 10619.6001
$ws.Range("M132").Value = -13673.4995
$ws.Range("N132").Value = -15679.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 76924220
$ws.Range("I22").Value = 111111500
$ws.Range("J22").Value = 2850.25
$ws.Range("K22").Value = 111111500
$ws.Range("L22").Value = 2850.25
$ws.Range("M22").Value = -111111205
$ws.Range("N22").Value = -3440.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 76924220
$ws.Range("I27").Value = 111111500
$ws.Range("J27").Value = 2850.25
$ws.Range("K27").Value = 111111500
$ws.Range("L27").Value = 2850.25
$ws.Range("M27").Value = -111111393
$ws.Range("N27").Value = -3064.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4322
$ws.Range("I40").Value = 4228
$ws.Range("K40").Value = 4228
$ws.Range("M40").Value = -4092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1760.8334
$ws.Range("I68").Value = 1020.9091
$ws.Range("K68").Value = 1020.9091
$ws.Range("M68").Value = -271.9091

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1760.8334
$ws.Range("I71").Value = 1020.9091
$ws.Range("K71").Value = 5104.5455
$ws.Range("M71").Value = -1360.5455

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3235.1667
$ws.Range("I82").Value = 2131.7144
$ws.Range("J82").Value = 4780
$ws.Range("K82").Value = 2131.7144
$ws.Range("L82").Value = 4780
$ws.Range("M82").Value = -1770.7144
$ws.Range("N82").Value = -5502

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3235.1667
$ws.Range("I85").Value = 2131.7144
$ws.Range("J85").Value = 4780
$ws.Range("K85").Value = 2131.7144
$ws.Range("L85").Value = 4780
$ws.Range("M85").Value = -883.7143999999998
$ws.Range("N85").Value = -7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H86").Value = 40000
$ws.Range("J86").Value = 40000
$ws.Range("L86").Value = 40000
$ws.Range("N86").Value = -42372

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H89").Value = 40000
$ws.Range("J89").Value = 40000
$ws.Range("L89").Value = 120000
$ws.Range("N89").Value = -131856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2042992.8
$ws.Range("I136").Value = 2858267.2
$ws.Range("J136").Value = 4806.2856
$ws.Range("K136").Value = 8574801.600000001
$ws.Range("L136").Value = 14418.8568
$ws.Range("M136").Value = -8572251.600000001
$ws.Range("N136").Value = -19518.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 34565
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 31333.334
$ws.Range("J76").Value = 31333.334
$ws.Range("L76").Value = 31333.334
$ws.Range("N76").Value = -31963.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 34565
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H79").Value = 31333.334
$ws.Range("J79").Value = 31333.334
$ws.Range("L79").Value = 31333.334
$ws.Range("N79").Value = -33517.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 33333.332
$ws.Range("J80").Value = 33333.332
$ws.Range("L80").Value = 33333.332
$ws.Range("N80").Value = -35329.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1748
$ws.Range("I81").Value = 1483.5
$ws.Range("J81").Value = 1880.25
$ws.Range("K81").Value = 2967
$ws.Range("L81").Value = 3760.5
$ws.Range("M81").Value = -1906
$ws.Range("N81").Value = -5882.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H83").Value = 33333.332
$ws.Range("J83").Value = 33333.332
$ws.Range("L83").Value = 99999.99600000001
$ws.Range("N83").Value = -109983.996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1748
$ws.Range("I84").Value = 1483.5
$ws.Range("J84").Value = 1880.25
$ws.Range("K84").Value = 14835
$ws.Range("L84").Value = 18802.5
$ws.Range("M84").Value = -9531
$ws.Range("N84").Value = -29410.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2859033
$ws.Range("I126").Value = 1245.96
$ws.Range("J126").Value = 10003501
$ws.Range("K126").Value = 3737.88
$ws.Range("L126").Value = 30010503
$ws.Range("M126").Value = -1267.88
$ws.Range("N126").Value = -30015443

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 112766.94
$ws.Range("I132").Value = 128911.336
$ws.Range("J132").Value = 34063
$ws.Range("K132").Value = 386734.008
$ws.Range("L132").Value = 102189
$ws.Range("M132").Value = -384204.008
$ws.Range("N132").Value = -107249
